## Quickdrive consists ZWNL 2022
## Adds a "Bkl" route-planning column (C) to the Standard sheet, inserts a
## blank spacer row above "Spots somda", moves the trailing comment down
## with it, and records the sort state used on the new column's data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Standard")

$excel.UserName = "Tom"

## --- 1. Insert a blank row above the old row 10 ("Spots somda"), pushing
##        everything below it down by one (old 10 -> 11, ... old 19 -> 20).
$ws.Rows.Item(10).Insert()

## --- 2. Re-home the "S-bord Amfpon stoptijd resetten" comment that used to
##        sit on B14; after the insert the matching data row is now 15.
$oldComment = $ws.Range("B14").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()
$ws.Range("B15").AddComment($commentText)

## --- 3. Record the sort that was applied to the new column's data block
##        while it is still empty, so reordering can't scramble real data.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("C3:C10"))
$sortObj.SetRange($ws.Range("C3:C10"))
$sortObj.Apply()

## --- 4. New column C: give it a width and fill in the route notes.
$ws.Columns.Item(3).ColumnWidth = 22.45

$ws.Range("C1").Value = "7322 Rhn-Bkl"

$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = "Done"

$ws.Range("B2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").WrapText = $true
$ws.Range("C3").Value = "Rhn-Db:`n- 7323`n- 3022*`n- 3125`n- 7325"

$ws.Range("B2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").WrapText = $true
$ws.Range("C4").Value = "Db-Ut:`n- 7425`n- 105`n- 3225`n- 6027`n- 3127`n- 3927 rijdt niet`n- 3924* rijdt niet`n- 7327"
$ws.Rows.Item(4).RowHeight = 128.25

$ws.Range("C3").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "Ut:`n- 527`n- 622`n- 2029`n- 3124*`n- 3527`n- 5725/4926`n- 5622/5627`n- 6024/8824`n- 6929"
$ws.Rows.Item(5).RowHeight = 142.5

$ws.Range("C3").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Ut-Bkl:`n- 7427`n- 3027`n- 3227`n- 829`n- 3129`n- 3124*"
$ws.Rows.Item(6).RowHeight = 99.75

$ws.Range("C3").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "Bkl: geen"

$ws.Range("B2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").NumberFormat = "@"

$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = "geen"

$ws.Range("B13").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = "- Klmos Mrg`n- Kalk/staaltrein Bnk 0907"
$ws.Rows.Item(12).RowHeight = 42.75

$ws.Range("B13").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "- Utoz`n- Ut noord"

$ws.Range("B2").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C15").PasteSpecial(-4122)

## --- 5. Update the view so it reflects where the edits were made.
$ws.Range("C12").Select()
$ws.Application.ActiveWindow.ScrollRow = 6
